$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report was inserted as row 48, pushing the previously
# existing rows 48-51 down to become rows 49-52 (their data is unchanged).
$ws.Rows.Item(48).Insert()

# Fill in the newly inserted row 48 with the latest week's data.
$ws.Range("A48").Value = 3
$ws.Range("B48").Value = "Femacal de La Calera"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44578
$ws.Range("E48").Value = 5
$ws.Range("F48").Value = 100112022
$ws.Range("G48").Value = "Arveja Verde"
$ws.Range("H48").Value = "Perfection"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 73
$ws.Range("K48").Value = 18000
$ws.Range("L48").Value = 19000
$ws.Range("M48").Value = 18521
$ws.Range("N48").Value = "`$/malla 25 kilos"
$ws.Range("O48").Value = "Provincia de Talca"
$ws.Range("P48").Value = 741
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
